# rename master/slave in pptx
#
# Applies:
#   - datetimeFigureOut placeholder text "2014/10/15" -> "2014/12/2"
#     (slide master + every slide layout)
#   - "Master" -> "Target" / "Slave" -> "Standby" state-diagram labels
#     (shapes nested inside a group on slide 1)
#   - "start" -> "go" transition label on the Master/Slave diagram
#     (a two-paragraph textbox: "start" / "slave" -> "go" / "slave")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" footer placeholder text wherever it
#    appears (slide master + all custom layouts).
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
      $tr = $shp.TextFrame.TextRange
      if ($tr.Text -eq "2014/10/15") {
        $tr.Text = "2014/12/2"
      }
    }
  }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
  $layout = $master.CustomLayouts.Item($li)
  Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1 contains a nested group ("グループ化 153") with the
#    Master / Slave roundRect labels and the start/slave text box.
#    Walk every shape with an explicit stack (NOT a recursive function -
#    the PS COM host chokes on recursive calls carrying live COM
#    collections) and patch the matching text runs in place so
#    unrelated formatting is left untouched.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$stack = New-Object System.Collections.ArrayList
[void]$stack.Add($slide.Shapes)

while ($stack.Count -gt 0) {
  $top = $stack.Count - 1
  $shapes = $stack[$top]
  $stack.RemoveAt($top)

  for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)

    if ($shp.Type -eq 6) {
      # msoGroup -> queue its children for traversal
      [void]$stack.Add($shp.GroupItems)
      continue
    }

    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }

    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text

    if ($text -eq "Master") {
      $tr.Text = "Target"
    } elseif ($text -eq "Slave") {
      $tr.Text = "Standby"
    } elseif ($text -eq ("start" + [char]13 + "slave")) {
      # Only the first paragraph ("start") changes to "go"; the
      # second paragraph ("slave") is left untouched.
      $tr.Text = "go" + [char]13 + "slave"
    }
  }
}
